$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set in_service (column E) to TRUE for rows 10 through 15
foreach ($r in 10..15) {
    $ws.Cells.Item($r, 5).Value = $true
}
